$d = $word.ActiveDocument

# Locate the "FECHA DE CIRUGIA:" label and collapse to its end so the
# new content lands right after the existing label text, inside the
# same paragraph (mirrors the "FECHA DE FIRMA DE CONSENTIMIENTO:" line
# elsewhere in the form, which already has a tab + fill-in field).
$r = $d.Content
$found = $r.Find.Execute("FECHA DE CIRUGÍA:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $r.Collapse(0)

    # Four tab stops, still bold like the label, to push the fill-in
    # area to the right.
    $r.InsertAfter("`t`t`t`t")
    $r.Collapse(0)

    # A long run of ellipsis characters used as a hand-fill blank for
    # the CIE10 diagnosis code(s). Track where this text starts so we
    # can drop the bold formatting it inherited from the label/tabs.
    $dotsStart = $r.Start
    $r.InsertAfter("………")
    $r.Collapse(0)
    $r.InsertAfter("…………………")
    $dotsEnd = $r.End

    $dotsRange = $d.Range($dotsStart, $dotsEnd)
    $dotsRange.Font.Bold = $false
}
